$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number need to be forced to Text format
# first, otherwise Excel auto-converts the numeric-looking string into a real number
# (the workbook stores these price columns as text). A multi-area Range string only
# applies formatting to its first area here, so loop cell-by-cell instead.
$textFormatCells = @("D4","D5","D6","D7","D8","D9","D10","D11","D12","D14","D15","D16","D18","D19","D20","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($ref in $textFormatCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated cell values row by row.
$ws.Range("D2").Value = '30.375.13'
$ws.Range("E2").Value = '  +0.49%  '
$ws.Range("D3").Value = '1.870.14'
$ws.Range("E3").Value = '  +0.54%  '
$ws.Range("D4").Value = '0.9995'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '234.90'
$ws.Range("E5").Value = '  -0.60%  '
$ws.Range("D6").Value = '0.9995'
$ws.Range("E6").Value = '  -0.05%  '
$ws.Range("D7").Value = '0.4697'
$ws.Range("E7").Value = '  +0.33%  '
$ws.Range("D8").Value = '0.2874'
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("D9").Value = '0.06580'
$ws.Range("E9").Value = '  +0.50%  '
$ws.Range("D10").Value = '21.71'
$ws.Range("E10").Value = '  -0.31%  '
$ws.Range("D11").Value = '0.07889'
$ws.Range("E11").Value = '  -0.50%  '
$ws.Range("D12").Value = '96.56'
$ws.Range("E12").Value = '  -1.16%  '
$ws.Range("D13").Value = '1.864.51'
$ws.Range("E13").Value = '  +0.11%  '
$ws.Range("D14").Value = '0.6938'
$ws.Range("E14").Value = '  +1.73%  '
$ws.Range("D15").Value = '5.109'
$ws.Range("E15").Value = '  -1.39%  '
$ws.Range("D16").Value = '268.83'
$ws.Range("E16").Value = '  +0.18%  '
$ws.Range("D17").Value = '30.300.16'
$ws.Range("E17").Value = '  +0.25%  '
$ws.Range("D18").Value = '14.02'
$ws.Range("E18").Value = '  +1.77%  '
$ws.Range("D19").Value = '0.000007683'
$ws.Range("E19").Value = '  +3.25%  '
$ws.Range("D20").Value = '1.002'
$ws.Range("E20").Value = '  +0.19%  '
$ws.Range("D21").Value = '2.111.84'
$ws.Range("E21").Value = '  -0.17%  '
$ws.Range("D22").Value = '0.9993'
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("D23").Value = '5.243'
$ws.Range("E23").Value = '  -1.58%  '
$ws.Range("D24").Value = '6.186'
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("D25").Value = '9.399'
$ws.Range("E25").Value = '  +1.84%  '
$ws.Range("D26").Value = '167.67'
$ws.Range("E26").Value = '  +0.33%  '
$ws.Range("D27").Value = '18.87'
$ws.Range("E27").Value = '  -0.31%  '
$ws.Range("D28").Value = '1.948'
$ws.Range("E28").Value = '  -0.67%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '1.362'
$ws.Range("E29").Value = '  -1.72%  '
$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").Value = '0.09891'
$ws.Range("E30").Value = '  +0.49%  '
$ws.Range("D31").Value = '4.384'
$ws.Range("E31").Value = '  -0.12%  '
$ws.Range("D32").Value = '1.460'
$ws.Range("E32").Value = '  -0.96%  '
$ws.Range("D33").Value = '4.074'
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("D34").Value = '0.04761'
$ws.Range("E34").Value = '  +1.02%  '
$ws.Range("D35").Value = '1.137'
$ws.Range("E35").Value = '  +0.09%  '
$ws.Range("D36").Value = '0.7039'
$ws.Range("E36").Value = '  +0.08%  '
$ws.Range("D37").Value = '2.722'
$ws.Range("E37").Value = '  +0.59%  '
$ws.Range("D38").Value = '0.01875'
$ws.Range("E38").Value = '  -0.25%  '
$ws.Range("D39").Value = '2.803'
$ws.Range("E39").Value = '  +7.12%  '
$ws.Range("D40").Value = '6.243'
$ws.Range("E40").Value = '  -0.13%  '
$ws.Range("D41").Value = '73.42'
$ws.Range("E41").Value = '  -1.64%  '
$ws.Range("D42").Value = '1.956'
$ws.Range("E42").Value = '  +0.62%  '
$ws.Range("D43").Value = '0.8435'
$ws.Range("E43").Value = '  -0.36%  '
$ws.Range("D44").Value = '0.4178'
$ws.Range("E44").Value = '  +0.05%  '
$ws.Range("D45").Value = '0.9999'
$ws.Range("E45").Value = '  +0.07%  '
$ws.Range("D46").Value = '102.71'
$ws.Range("E46").Value = '  -0.56%  '
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '967.12'
$ws.Range("E47").Value = '  +1.07%  '
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").Value = '7.135'
$ws.Range("E48").Value = '  -0.55%  '
$ws.Range("D49").Value = '9.109'
$ws.Range("E49").Value = '  -1.26%  '
$ws.Range("D50").Value = '34.55'
$ws.Range("E50").Value = '  +1.13%  '
$ws.Range("D51").Value = '0.05677'
$ws.Range("E51").Value = '  +0.25%  '
